$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("config")

# Resize the table "Tabla1" to include the new row 7 (A1:C7 instead of A1:C6)
$table = $ws.ListObjects.Item("Tabla1")
$table.Resize($ws.Range("A1:C7"))

# Update row 5: TO changes from BRL to GBP (FROM stays USD)
$ws.Range("B5").Value = "GBP"
$ws.Range("C5").Value = "USD"

# Update row 6: TO changes from JPY to MXN, FROM changes from COP to HNL
$ws.Range("B6").Value = "MXN"
$ws.Range("C6").Value = "HNL"

# Fill in row 7 (previously empty): TO=JPY, FROM=USD
$ws.Range("B7").Value = "JPY"
$ws.Range("C7").Value = "USD"
$ws.Range("A7").Formula = '=+CONCATENATE("https://www.oanda.com/currency-converter/es/?from=",Tabla1[[#This Row],[TO]],"&to=",Tabla1[[#This Row],[FROM]],"&amount=1")'

# Recalculate formulas so column A (URL) updates
$wb.Application.Calculate()

# Update the active selection to match the diff (B8 selected)
$ws.Range("B8").Select()

# Update the workbook view/window size to match the diff
$wb.Windows.Item(1).WindowState = -4143  # xlNormal, ensure not maximized
$wb.Windows.Item(1).Left = 15
$wb.Windows.Item(1).Top = 15
$wb.Windows.Item(1).Width = 20460
$wb.Windows.Item(1).Height = 11490
